$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update data rows 6-11: company code, cost center/location, depreciation area, depreciation key ---
$ws.Range("B6").Value = "1710"
$ws.Range("S6").Value = "YB_1702"
$ws.Range("AT6").Value = "32"
$ws.Range("AZ6").Value = "SUL2"

$ws.Range("B7").Value = "1710"
$ws.Range("N7").Value = "17101201"
$ws.Range("AT7").Value = "32"
$ws.Range("AZ7").Value = "SUL2"

$ws.Range("B8").Value = "1710"
$ws.Range("N8").Value = "17101201"
$ws.Range("AT8").Value = "32"
$ws.Range("AZ8").Value = "SUL2"

$ws.Range("B9").Value = "1710"
$ws.Range("S9").Value = "YB_1702"
$ws.Range("AT9").Value = "32"
$ws.Range("AZ9").Value = "SUL2"

$ws.Range("B10").Value = "1710"
$ws.Range("N10").Value = "17101201"
$ws.Range("AT10").Value = "32"
$ws.Range("AZ10").Value = "SUL2"

$ws.Range("B11").Value = "1710"
$ws.Range("S11").Value = "YB_1702"
$ws.Range("AT11").Value = "32"
$ws.Range("AZ11").Value = "SUL2"

# --- Update the sheet view: scroll so column AP is top-left, select AZ6:AZ11 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 42
$ws.Range("AZ6:AZ11").Select()
